$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update vak (subject) -> lesson-hour-pattern mapping for single/double hour lessons
# and the new facilities ("FA") row codes read from input.
# Row 15 (IN)
$ws.Range("B15").Value = "H1|1C|1C"
$ws.Range("C15").Value = "V1|1C|1C"
$ws.Range("D15").Value = "H2|0"
$ws.Range("E15").Value = "V2|0"
$ws.Range("F15").Value = "H3|0"
$ws.Range("G15").Value = "V3|0"

# Row 2 (NL)
$ws.Range("B2").Value = "H1|1ALG|1ALG"
$ws.Range("C2").Value = "V1|1ALG|1ALG"
$ws.Range("D2").Value = "H2|1ALG|1ALG"
$ws.Range("E2").Value = "V2|1ALG|1ALG"
$ws.Range("F2").Value = "H3|1ALG|1ALG"
$ws.Range("G2").Value = "V3|1ALG|1ALG"

# Row 3 (EN)
$ws.Range("B3").Value = "H1|1ALG|1ALG"
$ws.Range("C3").Value = "V1|1ALG|1ALG"
$ws.Range("D3").Value = "H2|1ALG|1ALG"
$ws.Range("E3").Value = "V2|1ALG|1ALG"
$ws.Range("F3").Value = "H3|1ALG|1ALG"
$ws.Range("G3").Value = "V3|1ALG|1ALG"

# Row 4 (WI)
$ws.Range("B4").Value = "H1|1ALG|1ALG"
$ws.Range("C4").Value = "V1|1ALG|1ALG"
$ws.Range("D4").Value = "H2|1ALG|1ALG"
$ws.Range("E4").Value = "V2|1ALG|1ALG"
$ws.Range("F4").Value = "H3|1ALG|1ALG"
$ws.Range("G4").Value = "V3|1ALG|1ALG"

# Row 5 (LA)
$ws.Range("B5").Value = "H1|0"
$ws.Range("C5").Value = "V1|1ALG|1ALG"
$ws.Range("D5").Value = "H2|0"
$ws.Range("E5").Value = "V2|1ALG|1ALG"
$ws.Range("F5").Value = "H3|0"
$ws.Range("G5").Value = "V3|1ALG|1ALG"

# Row 6 (GR)
$ws.Range("B6").Value = "H1|0"
$ws.Range("C6").Value = "V1|1ALG|1ALG"
$ws.Range("D6").Value = "H2|0"
$ws.Range("E6").Value = "V2|1ALG|1ALG"
$ws.Range("F6").Value = "H3|0"
$ws.Range("G6").Value = "V3|1ALG|1ALG"

# Row 7 (LO)
$ws.Range("B7").Value = "H1|2"
$ws.Range("C7").Value = "V1|2"
$ws.Range("D7").Value = "H2|2"
$ws.Range("E7").Value = "V2|2"
$ws.Range("F7").Value = "H3|2"
$ws.Range("G7").Value = "V3|2"

# Row 8 (GS)
$ws.Range("B8").Value = "H1|1ALG|1ALG"
$ws.Range("C8").Value = "V1|1ALG|1ALG"
$ws.Range("D8").Value = "H2|1ALG|1ALG"
$ws.Range("E8").Value = "V2|1ALG|1ALG"
$ws.Range("F8").Value = "H3|1ALG|1ALG"
$ws.Range("G8").Value = "V3|1ALG|1ALG"

# Row 9 (AK)
$ws.Range("B9").Value = "H1|1ALG|1ALG"
$ws.Range("C9").Value = "V1|1ALG|1ALG"
$ws.Range("D9").Value = "H2|1ALG|1ALG"
$ws.Range("E9").Value = "V2|1ALG|1ALG"
$ws.Range("F9").Value = "H3|1ALG|1ALG"
$ws.Range("G9").Value = "V3|1ALG|1ALG"

# Row 10 (HV)
$ws.Range("B10").Value = "H1|1ALG|1ALG"
$ws.Range("C10").Value = "V1|1ALG|1ALG"
$ws.Range("D10").Value = "H2|1ALG|1ALG"
$ws.Range("E10").Value = "V2|1ALG|1ALG"
$ws.Range("F10").Value = "H3|1ALG|1ALG"
$ws.Range("G10").Value = "V3|1ALG|1ALG"

# Row 11 (ML)
$ws.Range("B11").Value = "H1|0"
$ws.Range("C11").Value = "V1|0"
$ws.Range("D11").Value = "H2|0"
$ws.Range("E11").Value = "V2|0"
$ws.Range("F11").Value = "H3|1ALG|1ALG"
$ws.Range("G11").Value = "V3|1ALG|1ALG"

# Row 12 (BI)
$ws.Range("B12").Value = "H1|1ALG|1ALG"
$ws.Range("C12").Value = "V1|1ALG|1ALG"
$ws.Range("D12").Value = "H2|1ALG|1ALG"
$ws.Range("E12").Value = "V2|1ALG|1ALG"
$ws.Range("F12").Value = "H3|1ALG|1ALG"
$ws.Range("G12").Value = "V3|1ALG|1ALG"

# Row 13 (NA)
$ws.Range("B13").Value = "H1|0"
$ws.Range("C13").Value = "V1|0"
$ws.Range("D13").Value = "H2|1ALG|1ALG"
$ws.Range("E13").Value = "V2|1ALG|1ALG"
$ws.Range("F13").Value = "H3|1ALG|1ALG"
$ws.Range("G13").Value = "V3|1ALG|1ALG"

# Row 14 (SK)
$ws.Range("B14").Value = "H1|0"
$ws.Range("C14").Value = "V1|0"
$ws.Range("D14").Value = "H2|1ALG|1ALG"
$ws.Range("E14").Value = "V2|1ALG|1ALG"
$ws.Range("F14").Value = "H3|1ALG|1ALG"
$ws.Range("G14").Value = "V3|1ALG|1ALG"

# Row 16 (FA)
$ws.Range("B16").Value = "H1|1ALG|1ALG"
$ws.Range("C16").Value = "V1|1ALG|1ALG"
$ws.Range("D16").Value = "H2|1ALG|1ALG"
$ws.Range("E16").Value = "V2|1ALG|1ALG"
$ws.Range("F16").Value = "H3|1ALG|1ALG"
$ws.Range("G16").Value = "V3|1ALG|1ALG"

# Column widths after the longer "1ALG"/"1C" codes were entered (Excel auto-sized B:G)
$ws.Columns.Item(2).ColumnWidth = 13.1   # B -> stored width 14
$ws.Columns.Item(3).ColumnWidth = 13     # C -> stored width 13.8333..
$ws.Columns.Item(4).ColumnWidth = 13.1   # D -> stored width 14
$ws.Columns.Item(5).ColumnWidth = 13     # E -> stored width 13.8333..
$ws.Columns.Item(6).ColumnWidth = 13.1   # F -> stored width 14
$ws.Columns.Item(7).ColumnWidth = 13     # G -> stored width 13.8333..

# Selection ended on E19 after entering the data
$ws.Range("E19").Select() | Out-Null
